$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 7 new rows starting at row 12 (pushes the existing "view" block and
# everything below it down by 7 rows, e.g. old row 12 -> new row 19).
$ws.Rows("12:18").Insert()

# Rows 6-13 (column D) hold the list of Activity classes. 3 new entries were
# added (EditionExerciceActivity, EditionExercicePlaylistActivity and
# ListeSonsActivity) alongside the 5 pre-existing ones.
$ws.Range("D6").Value  = "ChronometreActivity"
$ws.Range("D7").Value  = "ListeSequencesActivity"
$ws.Range("D8").Value  = "EdititionSequenceActivity"
$ws.Range("D9").Value  = "EditionExerciceActivity"

# New "autre" (other) sub-category with a couple of related classes.
$ws.Range("C14").Value = "autre"
$ws.Range("D14").Value = "ChronometreActivity"
$ws.Range("D15").Value = "ChronoService"

$ws.Range("D13").Value = "ListeSonsActivity"
$ws.Range("D10").Value = "EditionExercicePlaylistActivity"
$ws.Range("D11").Value = "AjoutExerciceActivity"
$ws.Range("D12").Value = "AjoutSequenceActivity"

# Re-sort the activity list alphabetically (matches the workbook's
# sortState ref="D6:D13" / sortCondition ref="D6").
$ws.Range("D6:D13").Sort($ws.Range("D6"))

# Update the visible selection to match the new layout.
$ws.Range("D18").Select()
